$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '68.675.95'
$ws.Cells.Item(2, 5).Value = '  +2.36%  '
$ws.Cells.Item(3, 4).Value = '3.750.93'
$ws.Cells.Item(3, 5).Value = '  +1.92%  '
$ws.Cells.Item(4, 5).Value = '  +0.04%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '601.05'
$ws.Cells.Item(5, 5).Value = '  +1.71%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '168.81'
$ws.Cells.Item(6, 5).Value = '  +1.17%  '
$ws.Cells.Item(7, 4).Value = '3.746.36'
$ws.Cells.Item(7, 5).Value = '  +1.86%  '
$ws.Cells.Item(8, 5).Value = '  +0.03%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.534'
$ws.Cells.Item(9, 5).Value = '  +2.69%  '
$ws.Cells.Item(10, 5).Value = '  +2.15%  '
$ws.Cells.Item(11, 5).Value = '  +3.01%  '
$ws.Cells.Item(12, 5).Value = '  +0.93%  '
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '38.20'
$ws.Cells.Item(14, 5).Value = '  +2.78%  '
$ws.Cells.Item(15, 4).Value = '4.377.37'
$ws.Cells.Item(15, 5).Value = '  +1.89%  '
$ws.Cells.Item(16, 4).Value = '3.751.41'
$ws.Cells.Item(16, 5).Value = '  +1.87%  '
$ws.Cells.Item(17, 4).Value = '68.694.58'
$ws.Cells.Item(17, 5).Value = '  +2.30%  '
$ws.Cells.Item(18, 5).Value = '  +2.91%  '
$ws.Cells.Item(19, 5).Value = '  +0.76%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '17.07'
$ws.Cells.Item(20, 5).Value = '  +0.56%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '10.91'
$ws.Cells.Item(21, 5).Value = '  +19.82%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '494.66'
$ws.Cells.Item(22, 5).Value = '  +1.82%  '
$ws.Cells.Item(23, 5).Value = '  +1.58%  '
$ws.Cells.Item(24, 5).Value = '  +7.74%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '85.36'
$ws.Cells.Item(25, 5).Value = '  +0.71%  '
$ws.Cells.Item(26, 5).Value = '  +1.53%  '
$ws.Cells.Item(27, 5).Value = '  +2.39%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '10.20'
$ws.Cells.Item(28, 5).Value = '  +2.86%  '
$ws.Cells.Item(29, 5).Value = '  +0.48%  '
$ws.Cells.Item(30, 5).Value = '  +7.60%  '
$ws.Cells.Item(31, 5).Value = '  +2.64%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '7.93'
$ws.Cells.Item(32, 5).Value = '  +3.12%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '31.92'
$ws.Cells.Item(33, 5).Value = '  +0.65%  '
$ws.Cells.Item(34, 4).Value = '3.897.10'
$ws.Cells.Item(34, 5).Value = '  +2.02%  '
$ws.Cells.Item(35, 4).Value = '3.686.47'
$ws.Cells.Item(35, 5).Value = '  +1.90%  '
$ws.Cells.Item(36, 5).Value = '  +1.95%  '
$ws.Cells.Item(37, 5).Value = '  -0.05%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '1.01'
$ws.Cells.Item(38, 5).Value = '  +1.97%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '5.86'
$ws.Cells.Item(39, 5).Value = '  +2.42%  '
$ws.Cells.Item(40, 5).Value = '  +1.05%  '
$ws.Cells.Item(41, 5).Value = '  +1.10%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '440.47'
$ws.Cells.Item(42, 5).Value = '  -0.03%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '48.87'
$ws.Cells.Item(43, 5).Value = '  +0.38%  '
$ws.Cells.Item(44, 5).Value = '  +1.49%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '2.87'
$ws.Cells.Item(45, 5).Value = '  +3.89%  '
$ws.Cells.Item(46, 5).Value = '  +2.49%  '
$ws.Cells.Item(47, 5).Value = '  -0.02%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '40.29'
$ws.Cells.Item(48, 5).Value = '  +1.41%  '
$ws.Cells.Item(49, 2).Value = 'Monero'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '141.59'
$ws.Cells.Item(49, 5).Value = '  +0.62%  '
$ws.Cells.Item(50, 2).Value = 'Maker'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(50, 4).Value = '2.813.38'
$ws.Cells.Item(50, 5).Value = '  +1.79%  '
$ws.Cells.Item(51, 5).Value = '  +3.34%  '
